{"js": "// Remove the \"Ver no Jupiter Salvar em pdf Salvar em docx\" line, the\n// \"\u00a9 2020 . Contact: ...\" footer line, and the blank paragraph that\n// precedes them (right after the \"LOM3238: Projeto Integrado I (Requisito)\"\n// requirement line), while leaving the trailing blank paragraph and the\n// page-break paragraph untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst footerText = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\n\nlet footerIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === footerText) {\n    footerIndex = i;\n    break;\n  }\n}\n\nif (footerIndex !== -1) {\n  const blankBefore = paragraphs.items[footerIndex - 1];\n  const footerPara = paragraphs.items[footerIndex];\n  const copyrightPara = paragraphs.items[footerIndex + 1];\n\n  // Delete from the bottom up so earlier deletions don't shift the\n  // references we still need.\n  copyrightPara.delete();\n  footerPara.delete();\n  blankBefore.delete();\n\n  await context.sync();\n}\n", "ps1": "# Remove the \"Ver no Jupiter Salvar em pdf Salvar em docx\" line, the\n# \"\u00a9 2020 . Contact: ...\" footer line, and the blank paragraph that\n# precedes them (right after the \"LOM3238: Projeto Integrado I (Requisito)\"\n# requirement line), while leaving the trailing blank paragraph and the\n# page-break paragraph untouched.\n\n$d = $word.ActiveDocument\n\n$footerText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n\n$startPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $footerText) {\n        $startPara = $p\n        break\n    }\n}\n\nif ($startPara -ne $null) {\n    $blankBefore = $startPara.Previous()\n    $copyrightPara = $startPara.Next()\n\n    $rangeToDelete = $d.Range($blankBefore.Range.Start, $copyrightPara.Range.End)\n    $rangeToDelete.Delete()\n}\n"}
